# Edit script: adds "PO Forecast" sheet with forecast data,
# and renames the "Requested quantity" headers on the existing
# "Weekly Quantity" / "Monthly Trend" sheets.

$wb = $excel.ActiveWorkbook

# --- 1) Rename existing column headers ------------------------------------
$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the new "PO Forecast" sheet at the end of the workbook --------
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the look of the other two sheets: bold/bordered header style copied
# from "Weekly Quantity"!A1:B1, and the date-formatted style copied from
# "Weekly Quantity"!A2.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A62").PasteSpecial(-4122)

$r1 = @(45004.99999999999, 6, -2.198707557767376, 15.95599162875396)
$r2 = @(45011.99999999999, 6, -2.316796901213228, 16.05982781159582)
$r3 = @(45018.99999999999, 6, -2.90408352056706, 14.88750300530262)
$r4 = @(45025.99999999999, 7, -2.319017722314295, 15.93929630430099)
$r5 = @(45081.99999999999, 7, -2.488370401319052, 16.4385989284518)
$r6 = @(45088.99999999999, 7, -1.866169076568601, 15.83188848734275)
$r7 = @(45109.99999999999, 7, -2.083005480210313, 15.77648527635568)
$r8 = @(45116.99999999999, 7, -1.941752835644891, 16.11441037005392)
$r9 = @(45123.99999999999, 7, -1.70356789441796, 15.85425569190644)
$r10 = @(45130.99999999999, 7, -2.070678626461689, 16.08939399867687)
$r11 = @(45137.99999999999, 7, -1.235544319435909, 16.47873315039377)
$r12 = @(45144.99999999999, 7, -1.156580867570858, 16.08201214190159)
$r13 = @(45151.99999999999, 7, -2.032750924610213, 15.89345930669856)
$r14 = @(45158.99999999999, 7, -1.575172043725771, 16.4949892547714)
$r15 = @(45165.99999999999, 8, -1.363755570495368, 16.53319825635184)
$r16 = @(45172.99999999999, 8, -1.914514585588557, 16.58163460968835)
$r17 = @(45179.99999999999, 8, -1.593810896480707, 16.40603127610296)
$r18 = @(45186.99999999999, 8, -1.75774843108849, 16.78290486499952)
$r19 = @(45193.99999999999, 8, -1.039952951782832, 17.14176642056405)
$r20 = @(45200.99999999999, 8, -1.463393512652801, 16.27373884928173)
$r21 = @(45207.99999999999, 8, -1.500407102447341, 17.56888816474871)
$r22 = @(45221.99999999999, 8, -1.043070663926212, 16.54827323224994)
$r23 = @(45228.99999999999, 8, -0.3351607790343997, 17.33189454863059)
$r24 = @(45242.99999999999, 8, -0.8677246484884831, 17.46879196210667)
$r25 = @(45249.99999999999, 8, -1.052413274107335, 16.8949993939701)
$r26 = @(45256.99999999999, 8, -0.7675339674339375, 17.41421235029283)
$r27 = @(45263.99999999999, 8, -0.3315376172530036, 17.82178633895337)
$r28 = @(45270.99999999999, 8, -0.399725643066613, 17.63142907084159)
$r29 = @(45277.99999999999, 8, -0.7734250365210734, 16.70297143883357)
$r30 = @(45298.99999999999, 8, -0.1045012585765759, 17.79261233125149)
$r31 = @(45305.99999999999, 9, -0.7882826401541692, 18.06302082567518)
$r32 = @(45312.99999999999, 9, -0.1589789160041149, 18.39929267162538)
$r33 = @(45319.99999999999, 9, -0.4676477982960508, 17.34508025880785)
$r34 = @(45326.99999999999, 9, -0.2831869615061288, 18.20669462413562)
$r35 = @(45333.99999999999, 9, 0.01960400167133091, 17.98451839130605)
$r36 = @(45347.99999999999, 9, 0.002554329700039753, 17.42700657003207)
$r37 = @(45354.99999999999, 9, -0.484317875566342, 18.02079564800975)
$r38 = @(45361.99999999999, 9, 0.3282007543490409, 18.04089071444495)
$r39 = @(45459.99999999999, 10, 1.114132855367099, 18.68303902536004)
$r40 = @(45466.99999999999, 10, 0.5342488418470964, 18.40988268680006)
$r41 = @(45473.99999999999, 10, 0.337811874572967, 18.25649943619793)
$r42 = @(45480.99999999999, 10, 0.8959793426694471, 18.51822180603415)
$r43 = @(45501.99999999999, 10, 1.313305306721399, 18.84840726483647)
$r44 = @(45515.99999999999, 10, 0.5992115514455423, 18.04524264685327)
$r45 = @(45543.99999999999, 10, 0.6688765542614952, 19.38125628576071)
$r46 = @(45550.99999999999, 10, 1.110472170281729, 19.20139828681635)
$r47 = @(45557.99999999999, 10, 0.999807200383757, 19.69859040462827)
$r48 = @(45564.99999999999, 10, 1.29819934550588, 19.7136853716622)
$r49 = @(45571.99999999999, 10, 1.682256409937092, 19.74225901067347)
$r50 = @(45578.99999999999, 10, 1.241695490545222, 20.31851730538663)
$r51 = @(45585.99999999999, 11, 1.139756342911898, 20.07180501637379)
$r52 = @(45592.99999999999, 11, 1.243882745876185, 19.28125058482449)
$r53 = @(45606.99999999999, 11, 2.285954976470862, 19.63985904607358)
$r54 = @(45613.99999999999, 11, 1.055678095837736, 20.01421859756109)
$r55 = @(45620.99999999999, 11, 1.510440140195995, 20.04393115086272)
$r56 = @(45627.99999999999, 11, 1.938981179841789, 19.89096720404497)
$r57 = @(45634.99999999999, 11, 1.527560595674953, 19.10676711527387)
$r58 = @(45641.99999999999, 11, 2.063795460349306, 19.64323308413168)
$r59 = @(45648.99999999999, 11, 2.217541256671195, 19.75898220530893)
$r60 = @(45655.99999999999, 11, 1.871406018627869, 20.32803335230662)
$r61 = @(45662.99999999999, 11, 1.415372423206864, 19.56508437592143)
$poData = @($r1, $r2, $r3, $r4, $r5, $r6, $r7, $r8, $r9, $r10, $r11, $r12, $r13, $r14, $r15, $r16, $r17, $r18, $r19, $r20, $r21, $r22, $r23, $r24, $r25, $r26, $r27, $r28, $r29, $r30, $r31, $r32, $r33, $r34, $r35, $r36, $r37, $r38, $r39, $r40, $r41, $r42, $r43, $r44, $r45, $r46, $r47, $r48, $r49, $r50, $r51, $r52, $r53, $r54, $r55, $r56, $r57, $r58, $r59, $r60, $r61)

# --- 3) Fill in the forecast data rows -------------------------------------
$r = 2
foreach ($row in $poData) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

[void]$wsForecast.Range("A1").Select()
